$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.NumberFormat = '@'
$r.Value = '66.411.90'
$r.Style = 'Normal'
$r = $ws.Range('E2')
$r.NumberFormat = '@'
$r.Value = '  +1.32%  '
$r.Style = 'Normal'
$r = $ws.Range('D3')
$r.NumberFormat = '@'
$r.Value = '3.774.36'
$r.Style = 'Normal'
$r = $ws.Range('E3')
$r.NumberFormat = '@'
$r.Value = '  +1.65%  '
$r.Style = 'Normal'
$r = $ws.Range('E4')
$r.NumberFormat = '@'
$r.Value = '  -0.31%  '
$r.Style = 'Normal'
$r = $ws.Range('D5')
$r.NumberFormat = '@'
$r.Value = '420.08'
$r.Style = 'Normal'
$r = $ws.Range('E5')
$r.NumberFormat = '@'
$r.Value = '  +2.25%  '
$r.Style = 'Normal'
$r = $ws.Range('D6')
$r.NumberFormat = '@'
$r.Value = '128.28'
$r.Style = 'Normal'
$r = $ws.Range('E6')
$r.NumberFormat = '@'
$r.Value = '  -3.87%  '
$r.Style = 'Normal'
$r = $ws.Range('D7')
$r.NumberFormat = '@'
$r.Value = '3.774.47'
$r.Style = 'Normal'
$r = $ws.Range('E7')
$r.NumberFormat = '@'
$r.Value = '  +1.93%  '
$r.Style = 'Normal'
$r = $ws.Range('D8')
$r.NumberFormat = '@'
$r.Value = '0.600'
$r.Style = 'Normal'
$r = $ws.Range('E8')
$r.NumberFormat = '@'
$r.Value = '  -4.38%  '
$r.Style = 'Normal'
$r = $ws.Range('E9')
$r.NumberFormat = '@'
$r.Value = '  -0.16%  '
$r.Style = 'Normal'
$r = $ws.Range('D10')
$r.NumberFormat = '@'
$r.Value = '0.716'
$r.Style = 'Normal'
$r = $ws.Range('E10')
$r.NumberFormat = '@'
$r.Value = '  -2.89%  '
$r.Style = 'Normal'
$r = $ws.Range('D11')
$r.NumberFormat = '@'
$r.Value = '0.159'
$r.Style = 'Normal'
$r = $ws.Range('E11')
$r.NumberFormat = '@'
$r.Value = '  -4.19%  '
$r.Style = 'Normal'
$r = $ws.Range('D12')
$r.NumberFormat = '@'
$r.Value = '0.0000343'
$r.Style = 'Normal'
$r = $ws.Range('E12')
$r.NumberFormat = '@'
$r.Value = '  +0.84%  '
$r.Style = 'Normal'
$r = $ws.Range('D13')
$r.NumberFormat = '@'
$r.Value = '39.85'
$r.Style = 'Normal'
$r = $ws.Range('E13')
$r.NumberFormat = '@'
$r.Value = '  -6.13%  '
$r.Style = 'Normal'
$r = $ws.Range('D14')
$r.NumberFormat = '@'
$r.Value = '4.392.53'
$r.Style = 'Normal'
$r = $ws.Range('E14')
$r.NumberFormat = '@'
$r.Value = '  +1.92%  '
$r.Style = 'Normal'
$r = $ws.Range('D15')
$r.NumberFormat = '@'
$r.Value = '10.00'
$r.Style = 'Normal'
$r = $ws.Range('E15')
$r.NumberFormat = '@'
$r.Value = '  -0.11%  '
$r.Style = 'Normal'
$r = $ws.Range('D16')
$r.NumberFormat = '@'
$r.Value = '15.61'
$r.Style = 'Normal'
$r = $ws.Range('E16')
$r.NumberFormat = '@'
$r.Value = '  +20.26%  '
$r.Style = 'Normal'
$r = $ws.Range('E17')
$r.NumberFormat = '@'
$r.Value = '  -1.00%  '
$r.Style = 'Normal'
$r = $ws.Range('D18')
$r.NumberFormat = '@'
$r.Value = '3.771.07'
$r.Style = 'Normal'
$r = $ws.Range('E18')
$r.NumberFormat = '@'
$r.Value = '  +1.58%  '
$r.Style = 'Normal'
$r = $ws.Range('D19')
$r.NumberFormat = '@'
$r.Value = '19.41'
$r.Style = 'Normal'
$r = $ws.Range('E19')
$r.NumberFormat = '@'
$r.Value = '  -3.69%  '
$r.Style = 'Normal'
$r = $ws.Range('D20')
$r.NumberFormat = '@'
$r.Value = '66.564.40'
$r.Style = 'Normal'
$r = $ws.Range('E20')
$r.NumberFormat = '@'
$r.Value = '  +0.95%  '
$r.Style = 'Normal'
$r = $ws.Range('D21')
$r.NumberFormat = '@'
$r.Value = '1.07'
$r.Style = 'Normal'
$r = $ws.Range('E21')
$r.NumberFormat = '@'
$r.Value = '  -2.36%  '
$r.Style = 'Normal'
$r = $ws.Range('D22')
$r.NumberFormat = '@'
$r.Value = '401.46'
$r.Style = 'Normal'
$r = $ws.Range('E22')
$r.NumberFormat = '@'
$r.Value = '  -5.19%  '
$r.Style = 'Normal'
$r = $ws.Range('D23')
$r.NumberFormat = '@'
$r.Value = '14.11'
$r.Style = 'Normal'
$r = $ws.Range('E23')
$r.NumberFormat = '@'
$r.Value = '  -6.35%  '
$r.Style = 'Normal'
$r = $ws.Range('D24')
$r.NumberFormat = '@'
$r.Value = '83.14'
$r.Style = 'Normal'
$r = $ws.Range('E24')
$r.NumberFormat = '@'
$r.Value = '  -4.57%  '
$r.Style = 'Normal'
$r = $ws.Range('D25')
$r.NumberFormat = '@'
$r.Value = '2.98'
$r.Style = 'Normal'
$r = $ws.Range('E25')
$r.NumberFormat = '@'
$r.Value = '  -1.70%  '
$r.Style = 'Normal'
$r = $ws.Range('D26')
$r.NumberFormat = '@'
$r.Value = '36.78'
$r.Style = 'Normal'
$r = $ws.Range('E26')
$r.NumberFormat = '@'
$r.Value = '  +0.99%  '
$r.Style = 'Normal'
$r = $ws.Range('D27')
$r.NumberFormat = '@'
$r.Value = '5.54'
$r.Style = 'Normal'
$r = $ws.Range('E27')
$r.NumberFormat = '@'
$r.Value = '  +7.24%  '
$r.Style = 'Normal'
$r = $ws.Range('D28')
$r.NumberFormat = '@'
$r.Value = '3.18'
$r.Style = 'Normal'
$r = $ws.Range('E28')
$r.NumberFormat = '@'
$r.Value = '  -1.13%  '
$r.Style = 'Normal'
$r = $ws.Range('D29')
$r.NumberFormat = '@'
$r.Value = '9.32'
$r.Style = 'Normal'
$r = $ws.Range('E29')
$r.NumberFormat = '@'
$r.Value = '  -2.16%  '
$r.Style = 'Normal'
$r = $ws.Range('D30')
$r.NumberFormat = '@'
$r.Value = '718.38'
$r.Style = 'Normal'
$r = $ws.Range('E30')
$r.NumberFormat = '@'
$r.Value = '  +4.54%  '
$r.Style = 'Normal'
$r = $ws.Range('D31')
$r.NumberFormat = '@'
$r.Value = '8.48'
$r.Style = 'Normal'
$r = $ws.Range('E31')
$r.NumberFormat = '@'
$r.Value = '  +21.35%  '
$r.Style = 'Normal'
$r = $ws.Range('D32')
$r.NumberFormat = '@'
$r.Value = '2.75'
$r.Style = 'Normal'
$r = $ws.Range('E32')
$r.NumberFormat = '@'
$r.Value = '  +0.95%  '
$r.Style = 'Normal'
$r = $ws.Range('D33')
$r.NumberFormat = '@'
$r.Value = '12.29'
$r.Style = 'Normal'
$r = $ws.Range('E33')
$r.NumberFormat = '@'
$r.Value = '  -2.33%  '
$r.Style = 'Normal'
$r = $ws.Range('E34')
$r.NumberFormat = '@'
$r.Value = '  -0.24%  '
$r.Style = 'Normal'
$r = $ws.Range('D35')
$r.NumberFormat = '@'
$r.Value = '0.998'
$r.Style = 'Normal'
$r = $ws.Range('E35')
$r.NumberFormat = '@'
$r.Value = '  -0.17%  '
$r.Style = 'Normal'
$r = $ws.Range('E36')
$r.NumberFormat = '@'
$r.Value = '  -6.39%  '
$r.Style = 'Normal'
$r = $ws.Range('D37')
$r.NumberFormat = '@'
$r.Value = '38.17'
$r.Style = 'Normal'
$r = $ws.Range('E37')
$r.NumberFormat = '@'
$r.Value = '  -7.90%  '
$r.Style = 'Normal'
$r = $ws.Range('D38')
$r.NumberFormat = '@'
$r.Value = '54.75'
$r.Style = 'Normal'
$r = $ws.Range('E38')
$r.NumberFormat = '@'
$r.Value = '  -1.85%  '
$r.Style = 'Normal'
$r = $ws.Range('D39')
$r.NumberFormat = '@'
$r.Value = '0.0₃0752'
$r.Style = 'Normal'
$r = $ws.Range('E39')
$r.NumberFormat = '@'
$r.Value = '  +13.39%  '
$r.Style = 'Normal'
$r = $ws.Range('D40')
$r.NumberFormat = '@'
$r.Value = '4.95'
$r.Style = 'Normal'
$r = $ws.Range('E40')
$r.NumberFormat = '@'
$r.Value = '  +17.02%  '
$r.Style = 'Normal'
$r = $ws.Range('D41')
$r.NumberFormat = '@'
$r.Value = '0.0447'
$r.Style = 'Normal'
$r = $ws.Range('E41')
$r.NumberFormat = '@'
$r.Value = '  -5.34%  '
$r.Style = 'Normal'
$r = $ws.Range('D42')
$r.NumberFormat = '@'
$r.Value = '2.90'
$r.Style = 'Normal'
$r = $ws.Range('E42')
$r.NumberFormat = '@'
$r.Value = '  -2.23%  '
$r.Style = 'Normal'
$r = $ws.Range('E43')
$r.NumberFormat = '@'
$r.Value = '  +0.72%  '
$r.Style = 'Normal'
$r = $ws.Range('D44')
$r.NumberFormat = '@'
$r.Value = '0.133'
$r.Style = 'Normal'
$r = $ws.Range('E44')
$r.NumberFormat = '@'
$r.Value = '  -5.51%  '
$r.Style = 'Normal'
$r = $ws.Range('D45')
$r.NumberFormat = '@'
$r.Value = '3.30'
$r.Style = 'Normal'
$r = $ws.Range('E45')
$r.NumberFormat = '@'
$r.Value = '  -3.07%  '
$r.Style = 'Normal'
$r = $ws.Range('D46')
$r.NumberFormat = '@'
$r.Value = '142.77'
$r.Style = 'Normal'
$r = $ws.Range('E46')
$r.NumberFormat = '@'
$r.Value = '  -1.53%  '
$r.Style = 'Normal'
$r = $ws.Range('D47')
$r.NumberFormat = '@'
$r.Value = '3.04'
$r.Style = 'Normal'
$r = $ws.Range('E47')
$r.NumberFormat = '@'
$r.Value = '  -1.04%  '
$r.Style = 'Normal'
$r = $ws.Range('D48')
$r.NumberFormat = '@'
$r.Value = '2.01'
$r.Style = 'Normal'
$r = $ws.Range('E48')
$r.NumberFormat = '@'
$r.Value = '  -4.00%  '
$r.Style = 'Normal'
$r = $ws.Range('B49')
$r.NumberFormat = '@'
$r.Value = 'TheGraph'
$r.Style = 'Normal'
$r = $ws.Range('C49')
$r.NumberFormat = '@'
$r.Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$r.Style = 'Normal'
$r = $ws.Range('D49')
$r.NumberFormat = '@'
$r.Value = '0.303'
$r.Style = 'Normal'
$r = $ws.Range('E49')
$r.NumberFormat = '@'
$r.Value = '  +3.92%  '
$r.Style = 'Normal'
$r = $ws.Range('E50')
$r.NumberFormat = '@'
$r.Value = '  +0.50%  '
$r.Style = 'Normal'
$r = $ws.Range('B51')
$r.NumberFormat = '@'
$r.Value = 'EnergySwap'
$r.Style = 'Normal'
$r = $ws.Range('C51')
$r.NumberFormat = '@'
$r.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$r.Style = 'Normal'
$r = $ws.Range('D51')
$r.NumberFormat = '@'
$r.Value = '25.32'
$r.Style = 'Normal'
$r = $ws.Range('E51')
$r.NumberFormat = '@'
$r.Value = '  -7.05%  '
$r.Style = 'Normal'
